# Regenerate the "K" (strikeouts) column (column G) for each game log row
# using the recalculated K values (previously based on Strike# logic).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    3  = 0
    4  = 3
    5  = 1
    6  = 2
    7  = 0
    8  = 1
    9  = 4
    10 = 1
    11 = 4
    12 = 3
    13 = 1
    14 = 0
    15 = 2
    16 = 1
    17 = 3
    18 = 1
    19 = 4
    20 = 0
    21 = 6
    22 = 2
    23 = 4
    24 = 5
    25 = 4
    26 = 4
    27 = 7
    28 = 3
    29 = 8
    30 = 3
    31 = 4
    32 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
